# Sprint Logs update: mark workout history rows complete, add final workout-history
# log entry, and fill in the next few blank placeholder rows with the new
# Exercise Log / supervisor meeting entries for the "Update Exercise
# Functionality" sprint work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Flip "Not Complete" -> "Complete" for the already-logged Workout
# History rows (62-64) now that the work has actually finished. ---
$ws.Cells.Item(62, 7).Value = "Complete"
$ws.Cells.Item(63, 7).Value = "Complete"
$ws.Cells.Item(64, 7).Value = "Complete"

# --- Step 2: Row 65 was a half-filled placeholder row; fill it in as the final
# "Complete Workout History Screen" log entry. ---
$ws.Cells.Item(65, 1).Value = 45342
$ws.Cells.Item(65, 2).Value = "11:30-03:00"
$ws.Cells.Item(65, 6).Value = "Complete Workout History Screen"
$ws.Cells.Item(65, 7).Value = "Complete"

# --- Step 3: Rows 66-70 were blank placeholder rows (identical styling to every
# other not-yet-used row further down the sheet). Fill the next 3 of them in
# directly with the new Exercise Log sprint entries - no row insertion is
# needed since the rows were already blank and waiting to be used. ---
$ws.Rows(66).RowHeight = 43.2
$ws.Rows(67).RowHeight = 28.8
$ws.Rows(68).RowHeight = 28.8

# Row 66: Fixed GitHub issue / started Exercise Log screen
$ws.Cells.Item(66, 1).Value = 45343
$ws.Cells.Item(66, 2).Value = "08:30-12:30"
$ws.Cells.Item(66, 3).Value = "Home"
$ws.Cells.Item(66, 4).Value = "Application"
$ws.Cells.Item(66, 5).Value = "Ali Suhail"
$ws.Cells.Item(66, 6).Value = "Fix Github problem (Long Paths & Size Limiter)`nStart Exercise Log Screen`nImplement Delete Exercise Log Functionality"
$ws.Cells.Item(66, 6).WrapText = $true
$ws.Cells.Item(66, 7).Value = "Complete"
$ws.Cells.Item(66, 8).Value = "Fixed Github Problem`nImplemented Delete Exercise Log Functionality`nCompleted 40%  Exercise Log Screen"
$ws.Cells.Item(66, 8).WrapText = $true

# Row 67: Completed Exercise Log screen + Update functionality
$ws.Cells.Item(67, 1).Value = 45344
$ws.Cells.Item(67, 2).Value = "22:00-01:30"
$ws.Cells.Item(67, 3).Value = "Home"
$ws.Cells.Item(67, 4).Value = "Application"
$ws.Cells.Item(67, 5).Value = "Ali Suhail"
$ws.Cells.Item(67, 6).Value = "Complete Exercise Log Screen`nImplement Update Exercise Log Functionality"
$ws.Cells.Item(67, 6).WrapText = $true
$ws.Cells.Item(67, 7).Value = "Complete"
$ws.Cells.Item(67, 8).Value = "Implemented Exercise Log Screen`nImplemented Update Exercise Log Functionality"
$ws.Cells.Item(67, 8).WrapText = $true

# Row 68: Supervisor meeting entry (Online), re-using the meeting formatting from
# other "Online" supervisor-meeting rows, plus the new centered time-range cell
# style for column B (h:mm number format, bordered, centered).
$ws.Cells.Item(68, 1).Value = 45345

$b68 = $ws.Cells.Item(68, 2)
$b68.Value = "13:30-14:10"
$b68.NumberFormat = "h:mm"
$b68.HorizontalAlignment = -4108
$b68.VerticalAlignment = -4108
$b68.Borders.LineStyle = 1

$ws.Cells.Item(68, 3).Value = "Online"

$d68 = $ws.Cells.Item(68, 4)
$d68.Value = "Supervisor`nMeeting"
$d68.HorizontalAlignment = -4108
$d68.VerticalAlignment = -4108
$d68.WrapText = $true

$e68 = $ws.Cells.Item(68, 5)
$e68.Value = "Ali Suhail`nEman Qaddoumi"
$e68.HorizontalAlignment = -4108
$e68.VerticalAlignment = -4108
$e68.WrapText = $true

$f68 = $ws.Cells.Item(68, 6)
$f68.Value = "Discuss Project Progress"
$f68.VerticalAlignment = -4160
$f68.WrapText = $false

$ws.Cells.Item(68, 7).Value = "Complete"

$h68 = $ws.Cells.Item(68, 8)
$h68.Value = ""
$h68.VerticalAlignment = -4160
$h68.WrapText = $false

# --- Step 4: The next two blank placeholder rows (69, 70) get their recurring
# Home/Application/Ali Suhail/Not-Complete defaults plus upcoming dates, same as
# every other not-yet-logged row further down the sheet. ---
$ws.Cells.Item(69, 1).Value = 45346
$ws.Cells.Item(69, 3).Value = "Home"
$ws.Cells.Item(69, 4).Value = "Application"
$ws.Cells.Item(69, 5).Value = "Ali Suhail"
$ws.Cells.Item(69, 7).Value = "Not Complete"

$ws.Cells.Item(70, 1).Value = 45347
$ws.Cells.Item(70, 3).Value = "Home"
$ws.Cells.Item(70, 4).Value = "Application"
$ws.Cells.Item(70, 5).Value = "Ali Suhail"
$ws.Cells.Item(70, 7).Value = "Not Complete"

# --- Step 5: Restore view state (scroll position / zoom / active selection) to
# reflect where the author left off editing. ---
$ws.Range("G69").Select()
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 59
$win.ScrollColumn = 1
